$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.488.42'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.25'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.36'
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6298'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  -1.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2905'
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.89'
$ws.Range("E10").Value = '  +1.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07738'
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.843.74'
$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.011'
$ws.Range("E13").Value = '  +0.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6819'
$ws.Range("E14").Value = '  +0.57%  '

$ws.Range("E15").Value = '  -1.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.46'
$ws.Range("E16").Value = '  -0.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.297'
$ws.Range("E17").Value = '  +3.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.507.28'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.47'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("E20").Value = '  +0.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.516'
$ws.Range("E22").Value = '  +1.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.00'
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.508'
$ws.Range("E25").Value = '  +0.79%  '

$ws.Range("E26").Value = '  -2.33%  '

$ws.Range("E27").Value = '  -0.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06639'
$ws.Range("E28").Value = '  +17.16%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.457'
$ws.Range("E29").Value = '  +1.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.486'
$ws.Range("E30").Value = '  +0.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.094'
$ws.Range("E31").Value = '  -0.54%  '

$ws.Range("E32").Value = '  +1.31%  '

$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("E34").Value = '  -1.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6943'
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.563'
$ws.Range("E36").Value = '  -0.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01867'
$ws.Range("E37").Value = '  +2.47%  '

$ws.Range("E38").Value = '  +4.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.251.66'
$ws.Range("E39").Value = '  +1.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.780'
$ws.Range("E40").Value = '  +5.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9381'
$ws.Range("E41").Value = '  +3.98%  '

$ws.Range("E42").Value = '  +0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.014.73'
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.23'
$ws.Range("E44").Value = '  -0.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.18'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.086'
$ws.Range("E46").Value = '  -0.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.720'
$ws.Range("E47").Value = '  +2.97%  '

$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.977'
$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3933'
$ws.Range("E50").Value = '  -0.91%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05686'
$ws.Range("E51").Value = '  -0.23%  '
